# Applies the numeric corrections described in the commit diff
# (profit-calculation refresh across all 8 job sheets).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 267.8889
$ws.Range("I5").Value = 130.5
$ws.Range("J5").Value = 377.8
$ws.Range("K5").Value = 130.5
$ws.Range("L5").Value = 377.8
$ws.Range("M5").Value = -15.5
$ws.Range("N5").Value = -607.8
$ws.Range("H11").Value = 255.73334
$ws.Range("I11").Value = 255.73334
$ws.Range("K11").Value = 255.73334
$ws.Range("M11").Value = -115.73334
$ws.Range("H12").Value = 2654.8
$ws.Range("I12").Value = 1068.5
$ws.Range("J12").Value = 9000
$ws.Range("K12").Value = 1068.5
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = -898.5
$ws.Range("N12").Value = -9340
$ws.Range("H38").Value = 59.68421
$ws.Range("I38").Value = 59.68421
$ws.Range("K38").Value = 179.05263
$ws.Range("M38").Value = 192.94737
$ws.Range("H64").Value = 7449.75
$ws.Range("I64").Value = 3485.5715
$ws.Range("J64").Value = 12999.6
$ws.Range("K64").Value = 3485.5715
$ws.Range("L64").Value = 12999.6
$ws.Range("M64").Value = -3237.5715
$ws.Range("N64").Value = -13495.6
$ws.Range("H67").Value = 7449.75
$ws.Range("I67").Value = 3485.5715
$ws.Range("J67").Value = 12999.6
$ws.Range("K67").Value = 3485.5715
$ws.Range("L67").Value = 12999.6
$ws.Range("M67").Value = -2627.5715
$ws.Range("N67").Value = -14715.6
$ws.Range("H74").Value = 160087.14
$ws.Range("I74").Value = 257652.5
$ws.Range("K74").Value = 257652.5
$ws.Range("M74").Value = -256716.5
$ws.Range("H77").Value = 160087.14
$ws.Range("I77").Value = 257652.5
$ws.Range("K77").Value = 1288262.5
$ws.Range("M77").Value = -1283582.5
$ws.Range("H133").Value = 79489.8
$ws.Range("J133").Value = 79489.8
$ws.Range("L133").Value = 79489.8
$ws.Range("N133").Value = -89609.8
$ws.Range("H135").Value = 1340.7142
$ws.Range("I135").Value = 888.2727
$ws.Range("K135").Value = 7994.454299999999
$ws.Range("M135").Value = -5459.454299999999
$ws.Range("H138").Value = 26888.238
$ws.Range("I138").Value = 114418.555
$ws.Range("J138").Value = 5597.081
$ws.Range("K138").Value = 343255.665
$ws.Range("L138").Value = 16791.243
$ws.Range("M138").Value = -338115.665
$ws.Range("N138").Value = -27071.243

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 32261102
$ws.Range("I2").Value = 45457490
$ws.Range("J2").Value = 3267.889
$ws.Range("K2").Value = 45457490
$ws.Range("L2").Value = 3267.889
$ws.Range("M2").Value = -45457377
$ws.Range("N2").Value = -3493.889
$ws.Range("H32").Value = 5393.914
$ws.Range("I32").Value = 5182.4062
$ws.Range("K32").Value = 5182.4062
$ws.Range("M32").Value = -4895.4062
$ws.Range("H74").Value = 13488.866
$ws.Range("I74").Value = 1523.8
$ws.Range("K74").Value = 1523.8
$ws.Range("M74").Value = -649.8
$ws.Range("H77").Value = 13488.866
$ws.Range("I77").Value = 1523.8
$ws.Range("K77").Value = 7619
$ws.Range("M77").Value = -3251
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H116").Value = 32261102
$ws.Range("I116").Value = 45457490
$ws.Range("J116").Value = 3267.889
$ws.Range("K116").Value = 45457490
$ws.Range("L116").Value = 3267.889
$ws.Range("M116").Value = -45455196
$ws.Range("N116").Value = -7855.889

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 32261102
$ws.Range("I3").Value = 45457490
$ws.Range("J3").Value = 3267.889
$ws.Range("K3").Value = 45457490
$ws.Range("L3").Value = 3267.889
$ws.Range("M3").Value = -45457376
$ws.Range("N3").Value = -3495.889
$ws.Range("H99").Value = 8552.6
$ws.Range("I99").Value = 300
$ws.Range("K99").Value = 300
$ws.Range("M99").Value = 1198
$ws.Range("H103").Value = 26064.857
$ws.Range("J103").Value = 26064.857
$ws.Range("L103").Value = 26064.857
$ws.Range("N103").Value = -28408.857
$ws.Range("H105").Value = 1766.35
$ws.Range("I105").Value = 1410.5834
$ws.Range("K105").Value = 1410.5834
$ws.Range("M105").Value = 336.4166

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17339.152
$ws.Range("I31").Value = 24001.467
$ws.Range("K31").Value = 24001.467
$ws.Range("M31").Value = -23706.467
$ws.Range("H34").Value = 17339.152
$ws.Range("I34").Value = 24001.467
$ws.Range("K34").Value = 24001.467
$ws.Range("M34").Value = -23799.467
$ws.Range("H132").Value = 144123.64
$ws.Range("I132").Value = 188146.38
$ws.Range("J132").Value = 12055.444
$ws.Range("K132").Value = 564439.14
$ws.Range("L132").Value = 36166.33199999999
$ws.Range("M132").Value = -561909.14
$ws.Range("N132").Value = -41226.33199999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2213.8333
$ws.Range("I39").Value = 891.5
$ws.Range("K39").Value = 2674.5
$ws.Range("M39").Value = -2380.5
$ws.Range("H55").Value = 812.5
$ws.Range("J55").Value = 2125
$ws.Range("L55").Value = 6375
$ws.Range("N55").Value = -6729

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7603249
$ws.Range("I11").Value = 10857856
$ws.Range("J11").Value = 9166
$ws.Range("K11").Value = 10857856
$ws.Range("L11").Value = 9166
$ws.Range("M11").Value = -10857717
$ws.Range("N11").Value = -9444
$ws.Range("H24").Value = 22857.143
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 22857.143
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 22857.143
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -23203.143
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H127").Value = 76246.75
$ws.Range("J127").Value = 76246.75
$ws.Range("L127").Value = 76246.75
$ws.Range("N127").Value = -86166.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 15000
$ws.Range("I20").Value = 15000
$ws.Range("K20").Value = 15000
$ws.Range("M20").Value = -14774
$ws.Range("H98").Value = 33899
$ws.Range("J98").Value = 33899
$ws.Range("L98").Value = 33899
$ws.Range("N98").Value = -39889
$ws.Range("H122").Value = 6568.7334
$ws.Range("I122").Value = 5491.625
$ws.Range("K122").Value = 16474.875
$ws.Range("M122").Value = -14024.875
$ws.Range("H132").Value = 3346.2888
$ws.Range("I132").Value = 3097.5527
$ws.Range("K132").Value = 9292.658100000001
$ws.Range("M132").Value = -6762.658100000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 26342.857
$ws.Range("J42").Value = 26342.857
$ws.Range("L42").Value = 26342.857
$ws.Range("N42").Value = -27098.857
$ws.Range("H132").Value = 1273.7542
$ws.Range("I132").Value = 1139.1321
$ws.Range("K132").Value = 3417.3963
$ws.Range("M132").Value = -887.3963000000003

